$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''41.998.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '''2.209.81'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.41%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '''240.16'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.51%  '
$ws.Range('D6').Value = '''0.624'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D7').Value = '''72.93'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.14%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = '''0.602'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('D10').Value = '''42.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('D11').Value = '''0.0950'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('D12').Value = '''7.05'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.35%  '
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').Value = '''2.542.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').Value = '''14.13'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('D16').Value = '''0.835'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.74%  '
$ws.Range('D17').Value = '''2.210.69'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('D18').Value = '''41.877.34'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').Value = '''0.0000107'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.66%  '
$ws.Range('D20').Value = '''72.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').Value = '''6.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('D22').Value = '''10.08'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +14.24%  '
$ws.Range('D23').Value = '''228.74'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.07%  '
$ws.Range('E24').Value = '  -7.16%  '
$ws.Range('D25').Value = '''11.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.06%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').Value = '''2.25'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.09%  '
$ws.Range('E29').Value = '  +1.08%  '
$ws.Range('D30').Value = '''166.97'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.19%  '
$ws.Range('D31').Value = '''20.51'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('D32').Value = '''5.58'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.65%  '
$ws.Range('D33').Value = '''0.0787'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.68%  '
$ws.Range('D34').Value = '''0.124'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('D35').Value = '''28.64'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.80%  '
$ws.Range('D36').Value = '''0.109'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.27%  '
$ws.Range('D37').Value = '''4.23'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.90%  '
$ws.Range('D38').Value = '''0.0299'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.53%  '
$ws.Range('D39').Value = '''13.20'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('D40').Value = '''2.11'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('D41').Value = '''64.41'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.61%  '
$ws.Range('D42').Value = '''5.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.94%  '
$ws.Range('D43').Value = '''0.197'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').Value = '''8.68'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.62%  '
$ws.Range('D45').Value = '''103.52'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.89%  '
$ws.Range('E46').Value = '  -2.22%  '
$ws.Range('D47').Value = '''2.37'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.14%  '
$ws.Range('D48').Value = '''1.10'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.05%  '
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '''2.420.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.03%  '
